$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 45, pushing the existing rows 45..139 down to 46..140.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new weekly price-report entry.
$ws.Cells.Item(45, 1).Value = 4
$ws.Cells.Item(45, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(45, 3).Value = "Los Lagos"
$ws.Cells.Item(45, 4).Value = 44519
$ws.Cells.Item(45, 5).Value = 10
$ws.Cells.Item(45, 6).Value = 100112039
$ws.Cells.Item(45, 7).Value = "Ciboulette"
$ws.Cells.Item(45, 8).Value = "Sin especificar"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 240
$ws.Cells.Item(45, 11).Value = 2500
$ws.Cells.Item(45, 12).Value = 2500
$ws.Cells.Item(45, 13).Value = 2500
$ws.Cells.Item(45, 14).Value = "$/docena de atados"
$ws.Cells.Item(45, 15).Value = "Región Metropolitana"
$ws.Cells.Item(45, 16).Value = 833
$ws.Cells.Item(45, 17).Value = 3
$ws.Cells.Item(45, 18).Value = "Hortaliza"
